# Insert a new timestamp column right before the "nom" column (column EL,
# i.e. the 142nd column). Everything from EL onward (nom, url_produit)
# shifts one column to the right (EL->EM, EM->EN), matching the target
# dimension A1:EN206.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newColIndex = 142   # column EL

$ws.Columns.Item($newColIndex).Insert()

# Header row: new scrape timestamp for this column.
$ws.Cells.Item(1, $newColIndex).Value2 = "2026-02-03 11:26:05"

# Data rows: the new column carries forward the latest known price
# (copied from the column immediately to its left, EK) whenever that
# price exists; rows with no current price (EK blank) are left blank,
# same as they already are right after the column insert.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $priceCell = $ws.Cells.Item($r, $newColIndex - 1)
    $price = $priceCell.Value2
    if (-not [string]::IsNullOrEmpty($price)) {
        $ws.Cells.Item($r, $newColIndex).Value2 = $price
    }
}
